$d = $word.ActiveDocument

# The document's single section has two Pearson-logo pictures living in its
# footers and one BTec-logo picture living in its header. Word's
# HeaderFooter collections are indexed 1=Primary, 2=FirstPage, 3=EvenPages;
# pictures are identified by their (stable) alt text rather than by index
# since that is robust regardless of which header/footer slot holds them.

$sec = $d.Sections.Item(1)

# --- Footers: the Pearson logo pictures get renamed from image2.png to image1.png
for ($i = 1; $i -le 3; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}

# --- Header: the BTec logo picture gets renamed from image1.jpg to image2.jpg
for ($i = 1; $i -le 3; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }
}
